$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 154893.39
$ws.Range("I28").Value = 167592.83
$ws.Range("J28").Value = 2500
$ws.Range("K28").Value = 167592.83
$ws.Range("L28").Value = 2500
$ws.Range("M28").Value = -167107.83
$ws.Range("N28").Value = -3470
# Row 32
$ws.Range("H32").Value = 1000
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1000
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 1000
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -1652
# Row 33
$ws.Range("H33").Value = 359.35
$ws.Range("I33").Value = 204.94444
$ws.Range("J33").Value = 1749
$ws.Range("K33").Value = 204.94444
$ws.Range("L33").Value = 1749
$ws.Range("M33").Value = 24.05556000000001
$ws.Range("N33").Value = -2207
# Row 103
$ws.Range("H103").Value = 1671
$ws.Range("J103").Value = 1950
$ws.Range("L103").Value = 5850
$ws.Range("N103").Value = -7022
# Row 132
$ws.Range("H132").Value = 1331.1549
$ws.Range("I132").Value = 1320.2258
$ws.Range("K132").Value = 3960.6774
$ws.Range("M132").Value = -1430.6774
# Row 135
$ws.Range("H135").Value = 3704.923
$ws.Range("I135").Value = 1895.875
$ws.Range("K135").Value = 17062.875
$ws.Range("M135").Value = -14527.875
# Row 138
$ws.Range("H138").Value = 2519.28
$ws.Range("I138").Value = 1270.7858
$ws.Range("J138").Value = 2722.5232
$ws.Range("K138").Value = 3812.3574
$ws.Range("L138").Value = 8167.569600000001
$ws.Range("M138").Value = 1327.6426
$ws.Range("N138").Value = -18447.5696

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2422.0356
$ws.Range("I2").Value = 2440.7693
$ws.Range("J2").Value = 2405.8
$ws.Range("K2").Value = 2440.7693
$ws.Range("L2").Value = 2405.8
$ws.Range("M2").Value = -2327.7693
$ws.Range("N2").Value = -2631.8
# Row 32
$ws.Range("H32").Value = 32799478
$ws.Range("I32").Value = 38585268
$ws.Range("J32").Value = 8933097
$ws.Range("K32").Value = 38585268
$ws.Range("L32").Value = 8933097
$ws.Range("M32").Value = -38584981
$ws.Range("N32").Value = -8933671
# Row 45
$ws.Range("H45").Value = 2241.5
$ws.Range("I45").Value = 1174.75
$ws.Range("J45").Value = 4375
$ws.Range("K45").Value = 1174.75
$ws.Range("L45").Value = 4375
$ws.Range("M45").Value = -797.75
$ws.Range("N45").Value = -5129
# Row 61
$ws.Range("H61").Value = 3464.36
$ws.Range("I61").Value = 3269.923
$ws.Range("K61").Value = 3269.923
$ws.Range("M61").Value = -3057.923
# Row 110
$ws.Range("H110").Value = 2273.5
$ws.Range("J110").Value = 3500
$ws.Range("L110").Value = 3500
$ws.Range("N110").Value = -7590
# Row 116
$ws.Range("H116").Value = 2422.0356
$ws.Range("I116").Value = 2440.7693
$ws.Range("J116").Value = 2405.8
$ws.Range("K116").Value = 2440.7693
$ws.Range("L116").Value = 2405.8
$ws.Range("M116").Value = -146.7692999999999
$ws.Range("N116").Value = -6993.8
# Row 132
$ws.Range("H132").Value = 211488.83
$ws.Range("I132").Value = 272458.28
$ws.Range("J132").Value = 6409.8184
$ws.Range("K132").Value = 817374.8400000001
$ws.Range("L132").Value = 19229.4552
$ws.Range("M132").Value = -814844.8400000001
$ws.Range("N132").Value = -24289.4552
# Row 136
$ws.Range("H136").Value = 3464.36
$ws.Range("I136").Value = 3269.923
$ws.Range("K136").Value = 9809.769
$ws.Range("M136").Value = -7259.769

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2422.0356
$ws.Range("I3").Value = 2440.7693
$ws.Range("J3").Value = 2405.8
$ws.Range("K3").Value = 2440.7693
$ws.Range("L3").Value = 2405.8
$ws.Range("M3").Value = -2326.7693
$ws.Range("N3").Value = -2633.8
# Row 105
$ws.Range("H105").Value = 2450.9355
$ws.Range("I105").Value = 2131.9565
$ws.Range("J105").Value = 3368
$ws.Range("K105").Value = 2131.9565
$ws.Range("L105").Value = 3368
$ws.Range("M105").Value = -384.9564999999998
$ws.Range("N105").Value = -6862
# Row 134
$ws.Range("H134").Value = 2554230.5
$ws.Range("I134").Value = 3107900.5
$ws.Range("J134").Value = 7349.6
$ws.Range("K134").Value = 9323701.5
$ws.Range("L134").Value = 22048.8
$ws.Range("M134").Value = -9321166.5
$ws.Range("N134").Value = -27118.8

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2709.2156
$ws.Range("I58").Value = 2483.8
$ws.Range("J58").Value = 3528.9092
$ws.Range("K58").Value = 2483.8
$ws.Range("L58").Value = 3528.9092
$ws.Range("M58").Value = -2280.8
$ws.Range("N58").Value = -3934.9092
# Row 99
$ws.Range("H99").Value = 2711
$ws.Range("I99").Value = 2503.6667
$ws.Range("J99").Value = 3333
$ws.Range("K99").Value = 2503.6667
$ws.Range("L99").Value = 3333
$ws.Range("M99").Value = -1005.6667
$ws.Range("N99").Value = -6329
# Row 126
$ws.Range("H126").Value = 2711
$ws.Range("I126").Value = 2503.6667
$ws.Range("J126").Value = 3333
$ws.Range("K126").Value = 7511.000100000001
$ws.Range("L126").Value = 9999
$ws.Range("M126").Value = -5041.000100000001
$ws.Range("N126").Value = -14939
# Row 130
$ws.Range("H130").Value = 74994
$ws.Range("J130").Value = 74994
$ws.Range("L130").Value = 74994
$ws.Range("N130").Value = -85034
# Row 132
$ws.Range("H132").Value = 4074.6191
$ws.Range("I132").Value = 3987.1765
$ws.Range("K132").Value = 11961.5295
$ws.Range("M132").Value = -9431.529500000001
# Row 134
$ws.Range("H134").Value = 3065.8948
$ws.Range("I134").Value = 2956
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 8868
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -6333
$ws.Range("N134").Value = -17070
# Row 136
$ws.Range("H136").Value = 2709.2156
$ws.Range("I136").Value = 2483.8
$ws.Range("J136").Value = 3528.9092
$ws.Range("K136").Value = 7451.400000000001
$ws.Range("L136").Value = 10586.7276
$ws.Range("M136").Value = -4901.400000000001
$ws.Range("N136").Value = -15686.7276

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 2428.8
$ws.Range("I5").Value = 1715
$ws.Range("K5").Value = 5145
$ws.Range("M5").Value = -5033
# Row 135
$ws.Range("H135").Value = 2428.8
$ws.Range("I135").Value = 1715
$ws.Range("K135").Value = 15435
$ws.Range("M135").Value = -12900
# Row 140
$ws.Range("H140").Value = 2605.7368
$ws.Range("I140").Value = 864.4545000000001
$ws.Range("K140").Value = 2593.3635
$ws.Range("M140").Value = 2586.6365

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 2592.2173
$ws.Range("I132").Value = 1874.8422
$ws.Range("J132").Value = 5999.75
$ws.Range("K132").Value = 5624.5266
$ws.Range("L132").Value = 17999.25
$ws.Range("M132").Value = -3094.5266
$ws.Range("N132").Value = -23059.25

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3629.647
$ws.Range("I40").Value = 2800.2856
$ws.Range("K40").Value = 2800.2856
$ws.Range("M40").Value = -2664.2856
# Row 82
$ws.Range("H82").Value = 2893.6667
$ws.Range("I82").Value = 2870
$ws.Range("J82").Value = 2941
$ws.Range("K82").Value = 2870
$ws.Range("L82").Value = 2941
$ws.Range("M82").Value = -2509
$ws.Range("N82").Value = -3663
# Row 85
$ws.Range("H85").Value = 2893.6667
$ws.Range("I85").Value = 2870
$ws.Range("J85").Value = 2941
$ws.Range("K85").Value = 2870
$ws.Range("L85").Value = 2941
$ws.Range("M85").Value = -1622
$ws.Range("N85").Value = -5437
# Row 122
$ws.Range("H122").Value = 32300.3
$ws.Range("I122").Value = 34375.5
$ws.Range("J122").Value = 23999.5
$ws.Range("K122").Value = 103126.5
$ws.Range("L122").Value = 71998.5
$ws.Range("M122").Value = -100676.5
$ws.Range("N122").Value = -76898.5

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Range("H5").Value = 699999.5
$ws.Range("I5").Value = 599999
$ws.Range("K5").Value = 599999
$ws.Range("M5").Value = -599887
# Row 42
$ws.Range("H42").Value = 37345
$ws.Range("I42").Value = 29866.666
$ws.Range("J42").Value = 59780
$ws.Range("K42").Value = 29866.666
$ws.Range("L42").Value = 59780
$ws.Range("M42").Value = -29488.666
$ws.Range("N42").Value = -60536
# Row 113
$ws.Range("H113").Value = 795.125
$ws.Range("I113").Value = 765.8570999999999
$ws.Range("K113").Value = 2297.5713
$ws.Range("M113").Value = -127.5712999999996
# Row 122
$ws.Range("H122").Value = 8814.286
$ws.Range("I122").Value = 9949.5
$ws.Range("J122").Value = 2003
$ws.Range("K122").Value = 29848.5
$ws.Range("L122").Value = 6009
$ws.Range("M122").Value = -27398.5
$ws.Range("N122").Value = -10909
# Row 132
$ws.Range("H132").Value = 52955.75
$ws.Range("I132").Value = 55574.527
$ws.Range("J132").Value = 3199
$ws.Range("K132").Value = 166723.581
$ws.Range("L132").Value = 9597
$ws.Range("M132").Value = -164193.581
$ws.Range("N132").Value = -14657
# Row 136
$ws.Range("H136").Value = 3161.3
$ws.Range("I136").Value = 2122.8
$ws.Range("J136").Value = 4199.8
$ws.Range("K136").Value = 6368.400000000001
$ws.Range("L136").Value = 12599.4
$ws.Range("M136").Value = -3818.400000000001
$ws.Range("N136").Value = -17699.4
